$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the CE / QA "requested by" columns (F and G) for both data rows
# from the old user id "gsstjaya" to the new user id "gssprpra".
$ws.Range("F2").Value = "gssprpra"
$ws.Range("G2").Value = "gssprpra"
$ws.Range("F3").Value = "gssprpra"
$ws.Range("G3").Value = "gssprpra"

# Match the active selection left after the edit (F2:G3, active cell F2).
$ws.Range("F2:G3").Select()
